$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before the existing row 188, shifting all
# subsequent rows (old 188-309) down by one (to 189-310), matching the
# weekly refresh described in the commit message.
$ws.Rows.Item(188).Insert()

$ws.Range("A188").Value = 1
$ws.Range("B188").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C188").Value = "Arica y Parinacota"
$ws.Range("D188").Value = 44777
$ws.Range("E188").Value = 15
$ws.Range("F188").Value = 100114013
$ws.Range("G188").Value = "Zanahoria"
$ws.Range("H188").Value = "Sin especificar"
$ws.Range("I188").Value = "Primera"
$ws.Range("J188").Value = 70
$ws.Range("K188").Value = 21000
$ws.Range("L188").Value = 22000
$ws.Range("M188").Value = 21500
$ws.Range("N188").Value = "$/saco 25 kilos"
$ws.Range("O188").Value = "Región de Arica y Parinacota"
$ws.Range("P188").Value = 860
$ws.Range("Q188").Value = 25
$ws.Range("R188").Value = "Hortaliza"
